$wb = $excel.ActiveWorkbook

# --- StatOutput!A2:C2 -> new stat counts (0, 0, 3), kept as TEXT (shared-string) cells ---
$statOutput = $wb.Worksheets.Item("StatOutput")

foreach ($addr in @("A2", "B2", "C2")) {
    $statOutput.Range($addr).NumberFormat = "@"
}
$statOutput.Range("A2").Value = "0"
$statOutput.Range("B2").Value = "0"
$statOutput.Range("C2").Value = "3"
# Drop the temporary text number-format again so the cells fall back to the
# workbook's default (General) style, matching the original formatting.
$statOutput.Range("A2:C2").ClearFormats()

# --- StatOutput_Message!A18 -> updated Cypher query text (Akita -> Bullmastiff) ---
$statMsg = $wb.Worksheets.Item("StatOutput_Message")
$statMsg.Range("A18").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bullmastiff']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
